$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns for season record: Wins, Losses, Ties,
# matching the formatting (bold, centered, bordered) of the other headers.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every player row (2-47)
for ($row = 2; $row -le 47; $row++) {
    $ws.Cells.Item($row, 30).Value = 90
    $ws.Cells.Item($row, 31).Value = 72
    $ws.Cells.Item($row, 32).Value = 0
}
